$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crediti a inizi 2023")

# Insert a new row above the old row 5 ("type" / BS_CREDIT... row). This
# pushes that row and everything below it (incl. both tables) down by one.
$ws.Rows.Item(5).EntireRow.Insert()

# Populate the newly inserted Name/Value pair.
$ws.Range("A5").Value = "amount delta value"
$ws.Range("B5").Value = $false

# Materialize the (empty) C5/D5 cells like their neighbours above/below, and
# make sure none of the new row picks up the column's default "Comma" style.
$ws.Range("C5").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("B5:D5").Style = "Normal"

# Match the visible row height used by the other labelled rows (A3, A4, …).
$ws.Rows.Item(5).RowHeight = 15.5

# The inserted row falls inside "Table5" (A4:B6), so it grows to A4:B7.
# "CO__crediti2022" just shifts down (still 3 rows) to A16:G18.
$ws.ListObjects.Item("Table5").Resize($ws.Range("A4:B7"))
$ws.ListObjects.Item("CO__crediti2022").Resize($ws.Range("A16:G18"))

# Reflect the post-edit selection.
$ws.Range("A5").Select() | Out-Null
